# Auto-generated Excel COM-interop script
# Applies updated market-price values to the Leve profit tables across sheets
# per the scheduled runner's data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 725.8125
$ws.Range("I33").Value = 684.9167
$ws.Range("K33").Value = 684.9167
$ws.Range("M33").Value = -455.9167
$ws.Range("H52").Value = 298.9091
$ws.Range("J52").Value = 298.9091
$ws.Range("L52").Value = 896.7273
$ws.Range("N52").Value = -1216.7273
$ws.Range("H62").Value = 62544836
$ws.Range("I62").Value = 111121700
$ws.Range("J62").Value = 88877.14
$ws.Range("K62").Value = 111121700
$ws.Range("L62").Value = 88877.14
$ws.Range("M62").Value = -111121076
$ws.Range("N62").Value = -90125.14
$ws.Range("H65").Value = 62544836
$ws.Range("I65").Value = 111121700
$ws.Range("J65").Value = 88877.14
$ws.Range("K65").Value = 555608500
$ws.Range("L65").Value = 444385.7
$ws.Range("M65").Value = -555605380
$ws.Range("N65").Value = -450625.7
$ws.Range("H106").Value = 2701
$ws.Range("I106").Value = 2751.25
$ws.Range("K106").Value = 2751.25
$ws.Range("M106").Value = -2120.25
$ws.Range("H116").Value = 41673564
$ws.Range("I116").Value = 50004280
$ws.Range("K116").Value = 50004280
$ws.Range("M116").Value = -50000838
$ws.Range("H132").Value = 1267.8182
$ws.Range("I132").Value = 1293.1613
$ws.Range("K132").Value = 3879.4839
$ws.Range("M132").Value = -1349.4839
$ws.Range("H137").Value = 3587.35
$ws.Range("I137").Value = 4692
$ws.Range("K137").Value = 14076
$ws.Range("M137").Value = -11526
$ws.Range("H138").Value = 3883.62
$ws.Range("I138").Value = 1167.3572
$ws.Range("J138").Value = 7340.6816
$ws.Range("K138").Value = 3502.0716
$ws.Range("L138").Value = 22022.0448
$ws.Range("M138").Value = 1637.9284
$ws.Range("N138").Value = -32302.0448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4124.3335
$ws.Range("I2").Value = 2191.6667
$ws.Range("K2").Value = 2191.6667
$ws.Range("M2").Value = -2078.6667
$ws.Range("H61").Value = 6157.919
$ws.Range("I61").Value = 2519.1904
$ws.Range("K61").Value = 2519.1904
$ws.Range("M61").Value = -2307.1904
$ws.Range("H63").Value = 2326
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 2326
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
$ws.Range("H88").Value = 3999
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3999
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3999
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -4811
$ws.Range("H91").Value = 3999
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3999
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3999
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -6807
$ws.Range("H97").Value = 4638856.5
$ws.Range("I97").Value = 430
$ws.Range("K97").Value = 430
$ws.Range("M97").Value = 66
$ws.Range("H102").Value = 1092.3572
$ws.Range("I102").Value = 1028.5834
$ws.Range("K102").Value = 1028.5834
$ws.Range("M102").Value = 593.4166
$ws.Range("H116").Value = 4124.3335
$ws.Range("I116").Value = 2191.6667
$ws.Range("K116").Value = 2191.6667
$ws.Range("M116").Value = 102.3332999999998
$ws.Range("H132").Value = 5010.0713
$ws.Range("I132").Value = 2838.7856
$ws.Range("K132").Value = 8516.356800000001
$ws.Range("M132").Value = -5986.356800000001
$ws.Range("H136").Value = 6157.919
$ws.Range("I136").Value = 2519.1904
$ws.Range("K136").Value = 7557.5712
$ws.Range("M136").Value = -5007.5712

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4124.3335
$ws.Range("I3").Value = 2191.6667
$ws.Range("K3").Value = 2191.6667
$ws.Range("M3").Value = -2077.6667
$ws.Range("H20").Value = 9806017
$ws.Range("I20").Value = 23812186
$ws.Range("J20").Value = 1698.9
$ws.Range("K20").Value = 23812186
$ws.Range("L20").Value = 1698.9
$ws.Range("M20").Value = -23811939
$ws.Range("N20").Value = -2192.9
$ws.Range("H86").Value = 38503524
$ws.Range("I86").Value = 73922.14
$ws.Range("J86").Value = 83338060
$ws.Range("K86").Value = 73922.14
$ws.Range("L86").Value = 83338060
$ws.Range("M86").Value = -72799.14
$ws.Range("N86").Value = -83340306
$ws.Range("H89").Value = 38503524
$ws.Range("I89").Value = 73922.14
$ws.Range("J89").Value = 83338060
$ws.Range("K89").Value = 369610.7
$ws.Range("L89").Value = 416690300
$ws.Range("M89").Value = -363994.7
$ws.Range("N89").Value = -416701532
$ws.Range("H94").Value = 1455.0741
$ws.Range("I94").Value = 716.8125
$ws.Range("K94").Value = 716.8125
$ws.Range("M94").Value = -265.8125
$ws.Range("H99").Value = 2460046.2
$ws.Range("I99").Value = 2839.8333
$ws.Range("K99").Value = 2839.8333
$ws.Range("M99").Value = -1341.8333
$ws.Range("H105").Value = 3172.1292
$ws.Range("I105").Value = 2949.7827
$ws.Range("J105").Value = 3811.375
$ws.Range("K105").Value = 2949.7827
$ws.Range("L105").Value = 3811.375
$ws.Range("M105").Value = -1202.7827
$ws.Range("N105").Value = -7305.375
$ws.Range("H107").Value = 43271476
$ws.Range("I107").Value = 56250812
$ws.Range("J107").Value = 7016.1665
$ws.Range("K107").Value = 56250812
$ws.Range("L107").Value = 7016.1665
$ws.Range("M107").Value = -56248892
$ws.Range("N107").Value = -10856.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 68286.664
$ws.Range("J52").Value = 68286.664
$ws.Range("L52").Value = 68286.664
$ws.Range("N52").Value = -68874.664
$ws.Range("H86").Value = 5213446
$ws.Range("J86").Value = 4986
$ws.Range("L86").Value = 4986
$ws.Range("N86").Value = -7232
$ws.Range("H89").Value = 5213446
$ws.Range("J89").Value = 4986
$ws.Range("L89").Value = 24930
$ws.Range("N89").Value = -36162
$ws.Range("H107").Value = 1416.6666
$ws.Range("I107").Value = 1014.7619
$ws.Range("J107").Value = 2120
$ws.Range("K107").Value = 1014.7619
$ws.Range("L107").Value = 2120
$ws.Range("M107").Value = 905.2381
$ws.Range("N107").Value = -5960
$ws.Range("H132").Value = 8511.772000000001
$ws.Range("I132").Value = 3264.1428
$ws.Range("J132").Value = 10960.667
$ws.Range("K132").Value = 9792.428400000001
$ws.Range("L132").Value = 32882.001
$ws.Range("M132").Value = -7262.428400000001
$ws.Range("N132").Value = -37942.001
$ws.Range("H134").Value = 4833.06
$ws.Range("I134").Value = 1346.6923
$ws.Range("K134").Value = 4040.0769
$ws.Range("M134").Value = -1505.0769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 89992.336
$ws.Range("J52").Value = 89992.336
$ws.Range("L52").Value = 89992.336
$ws.Range("N52").Value = -90510.336
$ws.Range("H70").Value = 9477
$ws.Range("I70").Value = 8761.6
$ws.Range("K70").Value = 8761.6
$ws.Range("M70").Value = -8491.6
$ws.Range("H73").Value = 9477
$ws.Range("I73").Value = 8761.6
$ws.Range("K73").Value = 8761.6
$ws.Range("M73").Value = -7825.6
$ws.Range("H80").Value = 4298.375
$ws.Range("I80").Value = 2896
$ws.Range("J80").Value = 4498.7144
$ws.Range("K80").Value = 2896
$ws.Range("L80").Value = 4498.7144
$ws.Range("M80").Value = -1898
$ws.Range("N80").Value = -6494.7144
$ws.Range("H83").Value = 4298.375
$ws.Range("I83").Value = 2896
$ws.Range("J83").Value = 4498.7144
$ws.Range("K83").Value = 14480
$ws.Range("L83").Value = 22493.572
$ws.Range("M83").Value = -9488
$ws.Range("N83").Value = -32477.572
$ws.Range("H97").Value = 579.6
$ws.Range("I97").Value = 599.75
$ws.Range("J97").Value = 499
$ws.Range("K97").Value = 599.75
$ws.Range("L97").Value = 499
$ws.Range("M97").Value = -103.75
$ws.Range("N97").Value = -1491
$ws.Range("H102").Value = 3611
$ws.Range("I102").Value = 4531
$ws.Range("J102").Value = 1771
$ws.Range("K102").Value = 4531
$ws.Range("L102").Value = 1771
$ws.Range("M102").Value = -2909
$ws.Range("N102").Value = -5015
$ws.Range("H113").Value = 2938.25
$ws.Range("I113").Value = 2335.25
$ws.Range("K113").Value = 2335.25
$ws.Range("M113").Value = -165.25
$ws.Range("H132").Value = 9417.583000000001
$ws.Range("I132").Value = 2335.3333
$ws.Range("J132").Value = 16499.834
$ws.Range("K132").Value = 7005.999899999999
$ws.Range("L132").Value = 49499.50199999999
$ws.Range("M132").Value = -4475.999899999999
$ws.Range("N132").Value = -54559.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4992.7393
$ws.Range("I40").Value = 3987.7144
$ws.Range("J40").Value = 6556.1113
$ws.Range("K40").Value = 3987.7144
$ws.Range("L40").Value = 6556.1113
$ws.Range("M40").Value = -3851.7144
$ws.Range("N40").Value = -6828.1113
$ws.Range("H61").Value = 3244.4736
$ws.Range("I61").Value = 1496.9615
$ws.Range("K61").Value = 1496.9615
$ws.Range("M61").Value = -1294.9615
$ws.Range("H93").Value = 5654.32
$ws.Range("I93").Value = 6568.6
$ws.Range("J93").Value = 5044.8
$ws.Range("K93").Value = 6568.6
$ws.Range("L93").Value = 5044.8
$ws.Range("M93").Value = -5320.6
$ws.Range("N93").Value = -7540.8
$ws.Range("H113").Value = 3244.4736
$ws.Range("I113").Value = 1496.9615
$ws.Range("K113").Value = 1496.9615
$ws.Range("M113").Value = 673.0385000000001
